# The "FT fuel - Diesel" sheet lists, for each activity block, a set of
# technosphere/biosphere exchanges. For the "Syngas, RWGS, Production, for
# Fischer Tropsch process, hydrogen from coal gasification" activity, an
# extra "market group for electricity, low voltage" exchange (row 225) had
# snuck in alongside the Carbon monoxide + Hydrogen inputs, pushing the
# energy-allocation total above 1. Remove that stray exchange row; Excel
# shifts everything below it up by one row automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(225).Delete()

# Reflect the author's post-edit selection/scroll state.
$ws.Range("B218").Select()
